$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = "'29.118.86"
$ws.Range('E2').Value = "'  -0.29%  "

# Row 3
$ws.Range('D3').Value = "'1.841.49"
$ws.Range('E3').Value = "'  -0.49%  "

# Row 4
$ws.Range('D4').Value = "'0.9997"
$ws.Range('E4').Value = "'  +0.01%  "

# Row 5
$ws.Range('D5').Value = "'241.10"
$ws.Range('E5').Value = "'  -2.17%  "

# Row 6
$ws.Range('D6').Value = "'0.6864"
$ws.Range('E6').Value = "'  -1.75%  "

# Row 7
$ws.Range('E7').Value = "'  -0.01%  "

# Row 8
$ws.Range('D8').Value = "'0.3018"
$ws.Range('E8').Value = "'  -1.46%  "

# Row 9
$ws.Range('D9').Value = "'0.07462"
$ws.Range('E9').Value = "'  -3.39%  "

# Row 10
$ws.Range('D10').Value = "'23.13"
$ws.Range('E10').Value = "'  -1.74%  "

# Row 11
$ws.Range('D11').Value = "'0.07666"
$ws.Range('E11').Value = "'  -2.13%  "

# Row 12
$ws.Range('B12').Value = "Polkadot"
$ws.Range('C12').Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range('D12').Value = "'5.064"
$ws.Range('E12').Value = "'  -1.34%  "

# Row 13
$ws.Range('B13').Value = "WrappedEther"
$ws.Range('C13').Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range('D13').Value = "'1.831.12"
$ws.Range('E13').Value = "'  -0.97%  "

# Row 14
$ws.Range('D14').Value = "'0.6829"
$ws.Range('E14').Value = "'  -0.59%  "

# Row 15
$ws.Range('D15').Value = "'87.53"
$ws.Range('E15').Value = "'  -6.09%  "

# Row 16
$ws.Range('D16').Value = "'6.169"
$ws.Range('E16').Value = "'  -7.14%  "

# Row 17
$ws.Range('D17').Value = "'29.118.75"
$ws.Range('E17').Value = "'  -0.28%  "

# Row 18
$ws.Range('D18').Value = "'0.000008165"
$ws.Range('E18').Value = "'  -1.94%  "

# Row 19
$ws.Range('D19').Value = "'2.080.18"
$ws.Range('E19').Value = "'  -0.61%  "

# Row 20
$ws.Range('D20').Value = "'227.90"
$ws.Range('E20').Value = "'  -5.71%  "

# Row 21
$ws.Range('D21').Value = "'12.54"
$ws.Range('E21').Value = "'  -1.64%  "

# Row 22
$ws.Range('D22').Value = "'0.9999"
$ws.Range('E22').Value = "'  +0.02%  "

# Row 23
$ws.Range('D23').Value = "'7.394"
$ws.Range('E23').Value = "'  -1.80%  "

# Row 24
$ws.Range('E24').Value = "'  +0.00%  "

# Row 25
$ws.Range('D25').Value = "'0.1456"
$ws.Range('E25').Value = "'  -3.52%  "

# Row 26
$ws.Range('D26').Value = "'160.36"
$ws.Range('E26').Value = "'  +0.84%  "

# Row 27
$ws.Range('D27').Value = "'8.763"
$ws.Range('E27').Value = "'  -0.95%  "

# Row 28
$ws.Range('D28').Value = "'18.08"
$ws.Range('E28').Value = "'  -1.20%  "

# Row 29
$ws.Range('D29').Value = "'1.514"
$ws.Range('E29').Value = "'  -1.83%  "

# Row 30
$ws.Range('D30').Value = "'4.264"
$ws.Range('E30').Value = "'  +0.72%  "

# Row 31
$ws.Range('D31').Value = "'4.142"
$ws.Range('E31').Value = "'  -0.80%  "

# Row 32
$ws.Range('D32').Value = "'1.195"
$ws.Range('E32').Value = "'  +0.18%  "

# Row 33
$ws.Range('D33').Value = "'0.05192"
$ws.Range('E33').Value = "'  +1.37%  "

# Row 34
$ws.Range('D34').Value = "'0.7653"
$ws.Range('E34').Value = "'  -4.16%  "

# Row 35
$ws.Range('D35').Value = "'1.847"
$ws.Range('E35').Value = "'  -1.29%  "

# Row 36
$ws.Range('D36').Value = "'1.134"
$ws.Range('E36').Value = "'  -1.26%  "

# Row 37
$ws.Range('D37').Value = "'2.677"
$ws.Range('E37').Value = "'  -0.54%  "

# Row 38
$ws.Range('D38').Value = "'1.314.29"
$ws.Range('E38').Value = "'  +0.08%  "

# Row 39
$ws.Range('D39').Value = "'0.01835"
$ws.Range('E39').Value = "'  -1.99%  "

# Row 40
$ws.Range('E40').Value = "'  +0.46%  "

# Row 41
$ws.Range('D41').Value = "'0.9354"
$ws.Range('E41').Value = "'  -1.10%  "

# Row 42
$ws.Range('D42').Value = "'105.01"
$ws.Range('E42').Value = "'  -2.08%  "

# Row 43
$ws.Range('D43').Value = "'5.773"
$ws.Range('E43').Value = "'  -4.43%  "

# Row 44
$ws.Range('E44').Value = "'  -0.07%  "

# Row 45
$ws.Range('D45').Value = "'1.982.09"
$ws.Range('E45').Value = "'  -0.33%  "

# Row 46
$ws.Range('B46').Value = "Mantle"
$ws.Range('C46').Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range('D46').Value = "'0.5195"
$ws.Range('E46').Value = "'  +0.32%  "

# Row 47
$ws.Range('B47').Value = "BabyDogeCoin"
$ws.Range('C47').Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range('D47').Value = "'0.00000000123"
$ws.Range('E47').Value = "'  -0.45%  "

# Row 48
$ws.Range('D48').Value = "'64.88"
$ws.Range('E48').Value = "'  +1.19%  "

# Row 49
$ws.Range('D49').Value = "'9.556"
$ws.Range('E49').Value = "'  -1.79%  "

# Row 50
$ws.Range('D50').Value = "'1.773"
$ws.Range('E50').Value = "'  +0.46%  "

# Row 51
$ws.Range('E51').Value = "'  +0.91%  "
